$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 10: DAMSLTag (I) and DialogAct (J)
$ws.Range("I10").Value = "sv"
$ws.Range("J10").Value = "Statement-opinion"

# Row 18: DAMSLTag (I) and DialogAct (J)
$ws.Range("I18").Value = "ba"
$ws.Range("J18").Value = "Appreciation"

# Row 24: DAMSLTag (I) and DialogAct (J)
$ws.Range("I24").Value = "sv"
$ws.Range("J24").Value = "Statement-opinion"

# Row 31: DAMSLTag (I) and DialogAct (J)
$ws.Range("I31").Value = "sd"
$ws.Range("J31").Value = "Statement-non-opinion"
